$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.7 = 36365.22 pesos`n✅ 36365.22 pesos = 8.7 = 956.98 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 115
$wsTasas.Range("O10").Value = 4182
$wsTasas.Range("N12").Value = 4180
$wsTasas.Range("O12").Value = 110
